# Refresh the 合肥-漫展信息 exhibition listing on the "展览" and "全部类型"
# sheets to the upstream site state as of commit 456a3b4: the oldest
# event (old row 2, "肥东·原神&崩铁&崩坏only" on 2024-01-29) has already
# happened and drops off the list, every remaining event shifts up one
# row, a new event is appended at the end, and several "want to go" /
# price figures refresh to their latest live values.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Drop the old last data row (row 14) - the refreshed dataset has 12
    # events instead of 13, shrinking the used range to A1:I13.
    $ws.Rows.Item(14).Delete()

    # Row 2
    $ws.Range("A2").Value = 1
    $ws.Range("B2").NumberFormat = "@"
    $ws.Range("B2").Value = '2024-01-31'
    $ws.Range("B2").ClearFormats()
    $ws.Range("C2").Value = '肥西·原神&崩铁&崩坏only'
    $ws.Range("D2").Value = '仙满楼·麦肯希酒店 仙满楼·麦肯希酒店'
    $ws.Range("E2").Value = '2024.01.31 10:00-01.31 17:00'
    $ws.Range("F2").Value = 30
    $ws.Range("G2").Value = 55
    $ws.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=80944'
    $ws.Range("I2").Value = '//i0.hdslb.com/bfs/openplatform/202401/euD63Mlp1705479140627.jpeg'

    # Row 3
    $ws.Range("A3").Value = 2
    $ws.Range("B3").NumberFormat = "@"
    $ws.Range("B3").Value = '2024-02-03'
    $ws.Range("B3").ClearFormats()
    $ws.Range("C3").Value = '合肥·环形宇宙动漫游戏嘉年华'
    $ws.Range("D3").Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
    $ws.Range("E3").Value = '2024.02.03 09:30-02.04 17:00'
    $ws.Range("F3").Value = 6248
    $ws.Range("G3").Value = 65
    $ws.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=79963'
    $ws.Range("I3").Value = '//i0.hdslb.com/bfs/openplatform/202312/tBk3WVyX1702968658234.jpeg'

    # Row 4
    $ws.Range("A4").Value = 3
    $ws.Range("B4").NumberFormat = "@"
    $ws.Range("B4").Value = '2024-02-04'
    $ws.Range("B4").ClearFormats()
    $ws.Range("C4").Value = '合肥·环形宇宙动漫游戏嘉年华—吴晛专场'
    $ws.Range("D4").Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
    $ws.Range("E4").Value = '2024.02.04 11:30-02.04 17:00'
    $ws.Range("F4").Value = 174
    $ws.Range("G4").Value = 168
    $ws.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=80551'
    $ws.Range("I4").Value = '//i0.hdslb.com/bfs/openplatform/202401/MSS7qIQp1704695420767.jpeg'

    # Row 5
    $ws.Range("A5").Value = 4
    $ws.Range("B5").NumberFormat = "@"
    $ws.Range("B5").Value = '2024-02-04'
    $ws.Range("B5").ClearFormats()
    $ws.Range("C5").Value = '巢湖·原×铁×崩only'
    $ws.Range("D5").Value = '健康东路7号 巢湖国际饭店'
    $ws.Range("E5").Value = '2024.02.04 10:00-02.04 17:00'
    $ws.Range("F5").Value = 22
    $ws.Range("G5").Value = '不可售'
    $ws.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=80974'
    $ws.Range("I5").Value = '//i0.hdslb.com/bfs/openplatform/202401/wVVrdShB1705487994232.jpeg'

    # Row 6
    $ws.Range("A6").Value = 5
    $ws.Range("B6").NumberFormat = "@"
    $ws.Range("B6").Value = '2024-02-05'
    $ws.Range("B6").ClearFormats()
    $ws.Range("C6").Value = '合肥·六安lovelive only'
    $ws.Range("D6").Value = '经开区繁华大道与莲花路交叉口 百乐门大剧院'
    $ws.Range("E6").Value = '2024.02.05 09:00-02.05 17:00'
    $ws.Range("F6").Value = 37
    $ws.Range("G6").Value = 70
    $ws.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=81146'
    $ws.Range("I6").Value = '//i2.hdslb.com/bfs/openplatform/202401/QkgtYncY1705656564257.jpeg'

    # Row 7
    $ws.Range("A7").Value = 6
    $ws.Range("B7").NumberFormat = "@"
    $ws.Range("B7").Value = '2024-02-13'
    $ws.Range("B7").ClearFormats()
    $ws.Range("C7").Value = '合肥·新春AG动漫游戏盛典热血plus'
    $ws.Range("D7").Value = '山西路与太原路交叉口 挥动体育'
    $ws.Range("E7").Value = '2024.02.13 09:30-02.14 16:00'
    $ws.Range("F7").Value = 1883
    $ws.Range("G7").Value = 39.9
    $ws.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=80584'
    $ws.Range("I7").Value = '//i1.hdslb.com/bfs/openplatform/202401/yI94srFk1704703809648.jpeg'

    # Row 8
    $ws.Range("A8").Value = 7
    $ws.Range("B8").NumberFormat = "@"
    $ws.Range("B8").Value = '2024-02-17'
    $ws.Range("B8").ClearFormats()
    $ws.Range("C8").Value = '合肥·2024运动新春动漫庆典（全ip）'
    $ws.Range("D8").Value = '锦绣大道与清潭路交口东北角 李宁体育公园'
    $ws.Range("E8").Value = '2024.02.17 09:00-02.17 17:00'
    $ws.Range("F8").Value = 1412
    $ws.Range("G8").Value = 65
    $ws.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=79918'
    $ws.Range("I8").Value = '//i0.hdslb.com/bfs/openplatform/202312/vzuMc0sJ1702902061660.jpeg'

    # Row 9
    $ws.Range("A9").Value = 8
    $ws.Range("B9").NumberFormat = "@"
    $ws.Range("B9").Value = '2024-02-19'
    $ws.Range("B9").ClearFormats()
    $ws.Range("C9").Value = '合肥·安徽马娘only'
    $ws.Range("D9").Value = '桐城路与庐江路交叉口西南80米 赤阑桥文玩大厦'
    $ws.Range("E9").Value = '2024.02.19 09:00-02.19 17:00'
    $ws.Range("F9").Value = 293
    $ws.Range("G9").Value = 68
    $ws.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=78286'
    $ws.Range("I9").Value = '//i1.hdslb.com/bfs/openplatform/202311/721L5pIZ1699428443216.jpeg'

    # Row 10
    $ws.Range("A10").Value = 9
    $ws.Range("B10").NumberFormat = "@"
    $ws.Range("B10").Value = '2024-03-02'
    $ws.Range("B10").ClearFormats()
    $ws.Range("C10").Value = '合肥·星芒1.5动漫嘉年华'
    $ws.Range("D10").Value = '山西路与太原路交叉口 挥动体育'
    $ws.Range("E10").Value = '2024.03.02 09:30-03.02 17:30'
    $ws.Range("F10").Value = 954
    $ws.Range("G10").Value = 55
    $ws.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=81267'
    $ws.Range("I10").Value = '//i0.hdslb.com/bfs/openplatform/202401/GWidiefU1706003134747.jpeg'

    # Row 11
    $ws.Range("A11").Value = 10
    $ws.Range("B11").NumberFormat = "@"
    $ws.Range("B11").Value = '2024-03-16'
    $ws.Range("B11").ClearFormats()
    $ws.Range("C11").Value = '合肥·CW国潮动漫游戏嘉年华'
    $ws.Range("D11").Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
    $ws.Range("E11").Value = '2024.03.16 09:30-03.17 17:00'
    $ws.Range("F11").Value = 234
    $ws.Range("G11").Value = '不可售'
    $ws.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=81284'
    $ws.Range("I11").Value = '//i0.hdslb.com/bfs/openplatform/202401/38B92fWF1705995243803.jpeg'

    # Row 12
    $ws.Range("A12").Value = 11
    $ws.Range("B12").NumberFormat = "@"
    $ws.Range("B12").Value = '2024-04-04'
    $ws.Range("B12").ClearFormats()
    $ws.Range("C12").Value = '合肥· 第二届漫画城市动漫展 -故事再次开始'
    $ws.Range("D12").Value = '凤淮路与固镇路西北角 庐阳全民健身中心'
    $ws.Range("E12").Value = '2024.04.04 09:00-04.05 17:00'
    $ws.Range("F12").Value = 5577
    $ws.Range("G12").Value = 60
    $ws.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=78898'
    $ws.Range("I12").Value = '//i2.hdslb.com/bfs/openplatform/202311/244eBWip1700711342120.jpeg'

    # Row 13
    $ws.Range("A13").Value = 12
    $ws.Range("B13").NumberFormat = "@"
    $ws.Range("B13").Value = '2024-05-18'
    $ws.Range("B13").ClearFormats()
    $ws.Range("C13").Value = '合肥·梦时空SPO1动漫展'
    $ws.Range("D13").Value = '阜阳路16号 银瑞林国际大酒店'
    $ws.Range("E13").Value = '2024.05.18 10:00-05.18 17:00'
    $ws.Range("F13").Value = 70
    $ws.Range("G13").Value = 60
    $ws.Range("H13").Value = 'https://show.bilibili.com/platform/detail.html?id=80207'
    $ws.Range("I13").Value = '//i2.hdslb.com/bfs/openplatform/202312/tQQOHYE01703574162111.jpeg'

}
